# Commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
#
# Column G (header "K") holds a per-game stat that was regenerated for every
# data row (rows 2-74, i.e. A-index 0-72) in this save_data sheet. Row 70
# (A-index 68) already carried the freshly-regenerated value, so it is left
# untouched. Apply the newly computed "K" values cell-by-cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 2
    4  = 0
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 2
    11 = 2
    12 = 2
    13 = 0
    14 = 2
    15 = 0
    16 = 3
    17 = 1
    18 = 2
    19 = 2
    20 = 1
    21 = 2
    22 = 0
    23 = 0
    24 = 2
    25 = 1
    26 = 0
    27 = 2
    28 = 1
    29 = 1
    30 = 1
    31 = 0
    32 = 3
    33 = 2
    34 = 1
    35 = 2
    36 = 2
    37 = 1
    38 = 2
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 1
    45 = 0
    46 = 0
    47 = 1
    48 = 0
    49 = 2
    50 = 2
    51 = 1
    52 = 2
    53 = 0
    54 = 1
    55 = 2
    56 = 1
    57 = 1
    58 = 2
    59 = 2
    60 = 1
    61 = 1
    62 = 2
    63 = 2
    64 = 1
    65 = 0
    66 = 1
    67 = 2
    68 = 1
    69 = 1
    71 = 3
    72 = 1
    73 = 1
    74 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
